# Update cryptos list data (prices & 1h volume change) -- Mon Aug 12 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell (outside the used A1:E51 range) used to stage numeric-looking
# text so PasteSpecial(xlPasteValues) writes the literal string without
# leaving a quote-prefix style on the destination cell.
$scratch = $ws.Range('Z1')

$ws.Range('D2').Value = '58.538.02'
$ws.Range('E2').Value = '  -4.06%  '
$ws.Range('D3').Value = '2.539.22'
$ws.Range('E3').Value = '  -3.58%  '
$ws.Range('E4').Value = '  -0.03%  '
$scratch.Value = "'508.39"
$scratch.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  -4.16%  '
$scratch.Value = "'144.35"
$scratch.Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  -7.34%  '
$scratch.Value = "'0.998"
$scratch.Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -4.13%  '
$ws.Range('D9').Value = '2.544.32'
$ws.Range('E9').Value = '  -3.74%  '
$scratch.Value = "'6.16"
$scratch.Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  -7.74%  '
$ws.Range('E11').Value = '  -7.13%  '
$ws.Range('E12').Value = '  -5.32%  '
$ws.Range('E13').Value = '  -0.75%  '
$ws.Range('D14').Value = '2.984.02'
$ws.Range('E14').Value = '  -3.68%  '
$ws.Range('D15').Value = '58.516.39'
$ws.Range('E15').Value = '  -4.09%  '
$scratch.Value = "'20.76"
$scratch.Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  -5.87%  '
$ws.Range('E17').Value = '  -6.77%  '
$ws.Range('D18').Value = '2.543.31'
$ws.Range('E18').Value = '  -3.57%  '
$scratch.Value = "'4.54"
$scratch.Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  -5.12%  '
$scratch.Value = "'335.21"
$scratch.Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  -6.03%  '
$scratch.Value = "'10.10"
$scratch.Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  -5.17%  '
$ws.Range('E22').Value = '  +0.01%  '
$scratch.Value = "'5.96"
$scratch.Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  -4.66%  '
$scratch.Value = "'60.41"
$scratch.Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  -2.29%  '
$ws.Range('E25').Value = '  -4.93%  '
$scratch.Value = "'1.00"
$scratch.Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  +0.46%  '
$ws.Range('E27').Value = '  -5.52%  '
$ws.Range('D28').Value = '2.650.84'
$ws.Range('E28').Value = '  -3.70%  '
$ws.Range('D29').Value = '0.0₃0788'
$ws.Range('E29').Value = '  -9.67%  '
$scratch.Value = "'6.98"
$scratch.Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  -6.04%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$scratch.Value = "'5.86"
$scratch.Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  -4.46%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$scratch.Value = "'149.47"
$scratch.Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  -1.20%  '
$scratch.Value = "'18.56"
$scratch.Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  -5.01%  '
$scratch.Value = "'1.54"
$scratch.Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  -5.61%  '
$scratch.Value = "'0.920"
$scratch.Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  +4.01%  '
$scratch.Value = "'3.91"
$scratch.Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  -6.17%  '
$ws.Range('E38').Value = '  -7.58%  '
$scratch.Value = "'36.01"
$scratch.Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  -1.74%  '
$scratch.Value = "'0.826"
$scratch.Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  -12.03%  '
$ws.Range('E41').Value = '  -6.92%  '
$scratch.Value = "'284.21"
$scratch.Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  -5.31%  '
$ws.Range('E43').Value = '  -7.77%  '
$scratch.Value = "'0.0999"
$scratch.Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  -2.20%  '
$ws.Range('E45').Value = '  +0.02%  '
$scratch.Value = "'0.602"
$scratch.Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  -6.23%  '
$ws.Range('E47').Value = '  -5.02%  '
$scratch.Value = "'18.67"
$scratch.Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  -5.49%  '
$scratch.Value = "'10.29"
$scratch.Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  -0.57%  '
$scratch.Value = "'0.0227"
$scratch.Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  -5.23%  '
$ws.Range('E51').Value = '  -9.36%  '

$scratch.Clear()
